$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.306.51'
$ws.Range("E2").Value = '  -1.56%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.828.24'
$ws.Range("E3").Value = '  -1.22%  '

$ws.Range("E4").Value = '  -0.81%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.25'
$ws.Range("E5").Value = '  -1.81%  '

$ws.Range("E6").Value = '  -0.83%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4247'
$ws.Range("E7").Value = '  -1.92%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3715'
$ws.Range("E8").Value = '  -1.56%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07260'
$ws.Range("E9").Value = '  -1.73%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8660'
$ws.Range("E10").Value = '  -2.02%  '

$ws.Range("E11").Value = '  -2.33%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.828.37'
$ws.Range("E12").Value = '  -1.74%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.742'
$ws.Range("E13").Value = '  -0.17%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07097'
$ws.Range("E14").Value = '  -0.58%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.326'
$ws.Range("E15").Value = '  -2.73%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '89.89'
$ws.Range("E16").Value = '  +1.79%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.006'
$ws.Range("E17").Value = '  -1.05%  '

$ws.Range("E18").Value = '  -1.67%  '

$ws.Range("E19").Value = '  -0.67%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.11'
$ws.Range("E20").Value = '  -2.77%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.425.96'
$ws.Range("E21").Value = '  -1.13%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.135'
$ws.Range("E22").Value = '  -2.51%  '

$ws.Range("E23").Value = '  -2.59%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.062.18'
$ws.Range("E24").Value = '  -1.14%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.993'
$ws.Range("E25").Value = '  -1.82%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '152.85'
$ws.Range("E26").Value = '  -1.89%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.47'
$ws.Range("E27").Value = '  -0.87%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.181'
$ws.Range("E28").Value = '  +1.90%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.256'
$ws.Range("E29").Value = '  -2.97%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '116.64'
$ws.Range("E30").Value = '  -3.32%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08872'
$ws.Range("E31").Value = '  -1.10%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.200'
$ws.Range("E32").Value = '  -3.00%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7593'
$ws.Range("E33").Value = '  -2.36%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.466'
$ws.Range("E34").Value = '  -2.44%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.822'
$ws.Range("E35").Value = '  -3.53%  '

$ws.Range("E36").Value = '  -0.87%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.120'
$ws.Range("E37").Value = '  -2.36%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01980'
$ws.Range("E38").Value = '  +0.37%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05275'
$ws.Range("E39").Value = '  -1.12%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.348'
$ws.Range("E40").Value = '  +2.64%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.877'
$ws.Range("E41").Value = '  +0.55%  '

$ws.Range("E42").Value = '  +1.02%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5068'
$ws.Range("E43").Value = '  -2.40%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.700'
$ws.Range("E44").Value = '  -2.75%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.63'
$ws.Range("E45").Value = '  -1.37%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '107.71'
$ws.Range("E46").Value = '  -2.86%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4767'
$ws.Range("E47").Value = '  +0.47%  '

$ws.Range("E48").Value = '  -0.92%  '

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06395'
$ws.Range("E49").Value = '  -1.77%  '

$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.675'
$ws.Range("E50").Value = '  -2.41%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.861'
$ws.Range("E51").Value = '  -1.86%  '
